$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$tr.Text = $tr.Text + "`rPiemēram 2.9,2.1,1.0 utt."
